$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6718438267707825
$ws.Range("B1").Value = 1.836403608322144
$ws.Range("C1").Value = 3.471915721893311
$ws.Range("D1").Value = 2.86176586151123
$ws.Range("E1").Value = 1.65913724899292
